$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 9.108069666666667
$ws.Range("H2").Value = 27.324209
$ws.Range("I2").Value = 0.00155006418458712
$ws.Range("J2").Value = 0.00155006418458712
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 239.0839323333333
$ws.Range("N2").Value = 717.251797
$ws.Range("O2").Value = 0.4086975387666237
$ws.Range("P2").Value = 0.4086975387666237
$ws.Range("Q2").Value = 2177.593111872619
$ws.Range("R2").Value = 19598.33800685357
$ws.Range("S2").Value = 0.0006335074171710492
$ws.Range("T2").Value = 0.0006335074171710492

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 9.108069666666667
$ws.Range("H3").Value = 27.324209
$ws.Range("I3").Value = 0.00155006418458712
$ws.Range("J3").Value = 0.00155006418458712
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 117.0512696666667
$ws.Range("N3").Value = 351.153809
$ws.Range("O3").Value = 0.2000910950200451
$ws.Range("P3").Value = 0.2000910950200451
$ws.Range("Q3").Value = 1066.111118695787
$ws.Range("R3").Value = 9595.000068262081
$ws.Range("S3").Value = 0.00031015404004539
$ws.Range("T3").Value = 0.00031015404004539

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 9.108069666666667
$ws.Range("H4").Value = 27.324209
$ws.Range("I4").Value = 0.00155006418458712
$ws.Range("J4").Value = 0.00155006418458712
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 171.15883
$ws.Range("N4").Value = 513.47649
$ws.Range("O4").Value = 0.2925842480357353
$ws.Range("P4").Value = 0.2925842480357353
$ws.Range("Q4").Value = 1558.926547705157
$ws.Range("R4").Value = 14030.33892934641
$ws.Range("S4").Value = 0.0004535243638545476
$ws.Range("T4").Value = 0.0004535243638545476

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 9.108069666666667
$ws.Range("H5").Value = 27.324209
$ws.Range("I5").Value = 0.00155006418458712
$ws.Range("J5").Value = 0.00155006418458712
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 57.695868
$ws.Range("N5").Value = 173.087604
$ws.Range("O5").Value = 0.09862711817759588
$ws.Range("P5").Value = 0.09862711817759588
$ws.Range("Q5").Value = 525.497985222804
$ws.Range("R5").Value = 4729.481867005235
$ws.Range("S5").Value = 0.0001528783635161327
$ws.Range("T5").Value = 0.0001528783635161327

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 5771.873535333333
$ws.Range("H6").Value = 17315.620606
$ws.Range("I6").Value = 0.9822909543423312
$ws.Range("J6").Value = 0.9822909543423313
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 239.0839323333333
$ws.Range("N6").Value = 717.251797
$ws.Range("O6").Value = 0.4086975387666237
$ws.Range("P6").Value = 0.4086975387666237
$ws.Range("Q6").Value = 1379962.221758192
$ws.Range("R6").Value = 12419659.99582373
$ws.Range("S6").Value = 0.4014598953924287
$ws.Range("T6").Value = 0.4014598953924287

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 5771.873535333333
$ws.Range("H7").Value = 17315.620606
$ws.Range("I7").Value = 0.9822909543423312
$ws.Range("J7").Value = 0.9822909543423313
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 117.0512696666667
$ws.Range("N7").Value = 351.153809
$ws.Range("O7").Value = 0.2000910950200451
$ws.Range("P7").Value = 0.2000910950200451
$ws.Range("Q7").Value = 675605.1256661987
$ws.Range("R7").Value = 6080446.130995789
$ws.Range("S7").Value = 0.1965476726826421
$ws.Range("T7").Value = 0.1965476726826422

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 5771.873535333333
$ws.Range("H8").Value = 17315.620606
$ws.Range("I8").Value = 0.9822909543423312
$ws.Range("J8").Value = 0.9822909543423313
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 171.15883
$ws.Range("N8").Value = 513.47649
$ws.Range("O8").Value = 0.2925842480357353
$ws.Range("P8").Value = 0.2925842480357353
$ws.Range("Q8").Value = 987907.1212156169
$ws.Range("R8").Value = 8891164.090940554
$ws.Range("S8").Value = 0.2874028602285558
$ws.Range("T8").Value = 0.2874028602285558

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 5771.873535333333
$ws.Range("H9").Value = 17315.620606
$ws.Range("I9").Value = 0.9822909543423312
$ws.Range("J9").Value = 0.9822909543423313
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 57.695868
$ws.Range("N9").Value = 173.087604
$ws.Range("O9").Value = 0.09862711817759588
$ws.Range("P9").Value = 0.09862711817759588
$ws.Range("Q9").Value = 333013.2536072853
$ws.Range("R9").Value = 2997119.282465568
$ws.Range("S9").Value = 0.09688052603870453
$ws.Range("T9").Value = 0.09688052603870455

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 1.272029666666667
$ws.Range("H10").Value = 3.816089
$ws.Range("I10").Value = 0.0002164813950916887
$ws.Range("J10").Value = 0.0002164813950916887
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 239.0839323333333
$ws.Range("N10").Value = 717.251797
$ws.Range("O10").Value = 0.4086975387666237
$ws.Range("P10").Value = 0.4086975387666237
$ws.Range("Q10").Value = 304.1218547513259
$ws.Range("R10").Value = 2737.096692761933
$ws.Range("S10").Value = 0.00008847541336273821
$ws.Range("T10").Value = 0.00008847541336273823

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 1.272029666666667
$ws.Range("H11").Value = 3.816089
$ws.Range("I11").Value = 0.0002164813950916887
$ws.Range("J11").Value = 0.0002164813950916887
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 117.0512696666667
$ws.Range("N11").Value = 351.153809
$ws.Range("O11").Value = 0.2000910950200451
$ws.Range("P11").Value = 0.2000910950200451
$ws.Range("Q11").Value = 148.8926875370001
$ws.Range("R11").Value = 1340.034187833001
$ws.Range("S11").Value = 0.000043315999395363
$ws.Range("T11").Value = 0.00004331599939536301

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 1.272029666666667
$ws.Range("H12").Value = 3.816089
$ws.Range("I12").Value = 0.0002164813950916887
$ws.Range("J12").Value = 0.0002164813950916887
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 171.15883
$ws.Range("N12").Value = 513.47649
$ws.Range("O12").Value = 0.2925842480357353
$ws.Range("P12").Value = 0.2925842480357353
$ws.Range("Q12").Value = 217.7191094719566
$ws.Range("R12").Value = 1959.47198524761
$ws.Range("S12").Value = 0.00006333904619662866
$ws.Range("T12").Value = 0.00006333904619662866

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 1.272029666666667
$ws.Range("H13").Value = 3.816089
$ws.Range("I13").Value = 0.0002164813950916887
$ws.Range("J13").Value = 0.0002164813950916887
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 57.695868
$ws.Range("N13").Value = 173.087604
$ws.Range("O13").Value = 0.09862711817759588
$ws.Range("P13").Value = 0.09862711817759588
$ws.Range("Q13").Value = 73.39085574008399
$ws.Range("R13").Value = 660.517701660756
$ws.Range("S13").Value = 0.0000213509361369588
$ws.Range("T13").Value = 0.00002135093613695881

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 93.67702500000001
$ws.Range("H14").Value = 281.031075
$ws.Range("I14").Value = 0.01594250007799006
$ws.Range("J14").Value = 0.01594250007799006
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 239.0839323333333
$ws.Range("N14").Value = 717.251797
$ws.Range("O14").Value = 0.4086975387666237
$ws.Range("P14").Value = 0.4086975387666237
$ws.Range("Q14").Value = 22396.67150628798
$ws.Range("R14").Value = 201570.0435565918
$ws.Range("S14").Value = 0.006515660543661244
$ws.Range("T14").Value = 0.006515660543661244

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 93.67702500000001
$ws.Range("H15").Value = 281.031075
$ws.Range("I15").Value = 0.01594250007799006
$ws.Range("J15").Value = 0.01594250007799006
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 117.0512696666667
$ws.Range("N15").Value = 351.153809
$ws.Range("O15").Value = 0.2000910950200451
$ws.Range("P15").Value = 0.2000910950200451
$ws.Range("Q15").Value = 10965.01471484608
$ws.Range("R15").Value = 98685.1324336147
$ws.Range("S15").Value = 0.003189952297962185
$ws.Range("T15").Value = 0.003189952297962185

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 93.67702500000001
$ws.Range("H16").Value = 281.031075
$ws.Range("I16").Value = 0.01594250007799006
$ws.Range("J16").Value = 0.01594250007799006
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 171.15883
$ws.Range("N16").Value = 513.47649
$ws.Range("O16").Value = 0.2925842480357353
$ws.Range("P16").Value = 0.2925842480357353
$ws.Range("Q16").Value = 16033.64999688075
$ws.Range("R16").Value = 144302.8499719268
$ws.Range("S16").Value = 0.004664524397128373
$ws.Range("T16").Value = 0.004664524397128373

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 93.67702500000001
$ws.Range("H17").Value = 281.031075
$ws.Range("I17").Value = 0.01594250007799006
$ws.Range("J17").Value = 0.01594250007799006
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 57.695868
$ws.Range("N17").Value = 173.087604
$ws.Range("O17").Value = 0.09862711817759588
$ws.Range("P17").Value = 0.09862711817759588
$ws.Range("Q17").Value = 5404.777269032701
$ws.Range("R17").Value = 48642.9954212943
$ws.Range("S17").Value = 0.001572362839238257
$ws.Range("T17").Value = 0.001572362839238257

